$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 94, shifting existing rows 94:108 down to 95:109
$ws.Rows.Item(94).Insert()

# Copy the date number format (style) from the row below (now row 95, column D) onto the new D94 cell
$ws.Range("D95").Copy()
$ws.Range("D94").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# Populate the new row 94 with the weekly price entry data
$ws.Cells.Item(94, 1).Value = 10
$ws.Cells.Item(94, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(94, 3).Value = "La Araucanía"
$ws.Cells.Item(94, 4).Value = 44776
$ws.Cells.Item(94, 5).Value = 9
$ws.Cells.Item(94, 6).Value = 100114002
$ws.Cells.Item(94, 7).Value = "Camote"
$ws.Cells.Item(94, 8).Value = "Sin especificar"
$ws.Cells.Item(94, 9).Value = "Primera"
$ws.Cells.Item(94, 10).Value = 30
$ws.Cells.Item(94, 11).Value = 20000
$ws.Cells.Item(94, 12).Value = 20000
$ws.Cells.Item(94, 13).Value = 20000
$ws.Cells.Item(94, 14).Value = "$/malla 20 kilos"
$ws.Cells.Item(94, 15).Value = "Perú"
$ws.Cells.Item(94, 16).Value = 1000
$ws.Cells.Item(94, 17).Value = 20
$ws.Cells.Item(94, 18).Value = "Hortaliza"
